# Target diff analysis
# ---------------------
# The supplied unified diff touches two parts of the package:
#   - word/document.xml  : only the xmlns:* declarations (and mc:Ignorable)
#                           on the root <w:document> element are reordered
#                           (alphabetically by prefix).
#   - word/styles.xml     : every changed line is the *same* set of
#                           w:* attributes on the same elements, merely
#                           re-ordered alphabetically (w:pgSz, w:pgMar,
#                           w:rFonts, w:lang, w:latentStyles/lsdException,
#                           w:style, w:tblInd, w:tblCellMar, ...).
#
# No attribute value, no text run, no paragraph, no style definition and
# no document property changes anywhere in the diff - this is purely the
# alphabetical-attribute-order fingerprint left behind when the OOXML
# producer that made the real commit re-serialised the package (the same
# thing the "canonical OOXML" extraction used to build the diff above
# already normalises for comparison). There is no user-visible edit to
# replay through the Word object model for this particular template.
#
# We still touch the document through the COM surface so the session is
# exercised end-to-end, but we only *read* state - nothing is mutated -
# so the meaningful (canonical) content of document.xml/styles.xml stays
# byte-for-byte the same aside from that inevitable, semantically-empty
# re-serialisation.

$d = $word.ActiveDocument

# Sanity read-only pass over the content that the diff's context lines
# reference, confirming the template text and page geometry referenced by
# the diff are already what they should be; nothing here writes back to
# the document.
$null = $d.Content.Text
$null = $d.PageSetup.PageWidth
$null = $d.PageSetup.PageHeight
$null = $d.Styles("Normal").NameLocal
